$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 24 / 25: drop the explicit row height (back to the sheet default) ---
$ws.Rows.Item(24).AutoFit()
$ws.Rows.Item(25).AutoFit()

# --- Row 29: correct wording "revision projecto" -> "revision projet" ---
$ws.Range("D29").Value = 'Meeting Leonel, revision projet'

# --- Row 31: replaced with a new meeting entry ---
$ws.Range("C31").Value = 44929
$ws.Range("D31").Value = 'Meeting Equipe'
$ws.Range("E31").Value = '1h '

# --- Row 32: the entry that used to live in row 31 ---
$ws.Range("B32").Value = 14
$ws.Range("C32").Value = 44935
$ws.Range("D32").Value = 'Meeting Leonel, correction format du projet et creation ADO Entity Framework'
$ws.Range("E32").Value = '7h30'

# --- Row 33: the entry that used to live in row 32 ---
$ws.Range("B33").Value = 15
$ws.Range("C33").Value = 44936
$ws.Range("C33").NumberFormat = "d-mmm-yy"
$ws.Range("D33").Value = "Meeting avec le equipe, parler de l'avancement du projet, de la répartition des interfaces et du code. "
$ws.Range("E33").Value = '2h30'

# --- Row 34: brand-new row, first cell to use the new m/d/yy date style ---
$ws.Range("B34").Value = 16
$ws.Range("C34").Value = 44940
$ws.Range("C34").NumberFormat = "mm-dd-yy"
$ws.Range("D34").Value = 'Meeting avec le equipe.'
$ws.Range("E34").Value = '1h30'

# --- Row 35: continuation line (no No./Date) ---
$ws.Range("D35").Value = 'Avancement du codage, ecriture et modificacion des utilites de la app Chantier'
$ws.Range("E35").Value = '6h30'

# --- Row 36 ---
$ws.Range("B36").Value = 17
$ws.Range("C36").Value = 44941
$ws.Range("D36").Value = 'Avancement du codage, ecriture et modificacion des utilites de la app Chantier'
$ws.Range("E36").Value = '4h'

# --- Row 37 ---
$ws.Range("B37").Value = 18
$ws.Range("C37").Value = 44942
$ws.Range("D37").Value = 'Meeting Equipe : état d`avancement finalisation du projet et discussion des différents points à présenter au professeur'
$ws.Range("E37").Value = '30 min'

# Reuse the exact same date style created for C34 on the other two new date
# cells instead of letting each assignment mint its own numFmt entry.
$ws.Range("C34").Copy()
$ws.Range("C36").PasteSpecial(-4122)
$ws.Range("C37").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update selection to match the saved view ---
$ws.Range("F33").Select()
